$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.044.24"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.612.34"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.63"
$ws.Range("E5").Value = "  +10.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "568.12"
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("D7").Value = "3.608.43"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.674"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.02"
$ws.Range("E11").Value = "  +14.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.152"
$ws.Range("E12").Value = "  +3.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000288"
$ws.Range("E13").Value = "  +11.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.03"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("D15").Value = "4.192.71"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "3.616.61"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").Value = "67.885.74"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.40"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.99"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.01"
$ws.Range("E23").Value = "  +15.62%  "
$ws.Range("E24").Value = "  -4.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.50"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.96"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.61"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.94"
$ws.Range("E28").Value = "  +11.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.12"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.45"
$ws.Range("E30").Value = "  +20.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.47"
$ws.Range("E31").Value = "  +5.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.67"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "678.55"
$ws.Range("E33").Value = "  +9.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.23"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.93"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.27"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.425"
$ws.Range("E38").Value = "  +8.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "0.0₃0772"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  +15.41%  "
$ws.Range("D42").Value = "3.252.74"
$ws.Range("E42").Value = "  +7.04%  "
$ws.Range("E43").Value = "  +4.53%  "
$ws.Range("E44").Value = "  +11.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.05"
$ws.Range("E45").Value = "  +30.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0420"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  +10.38%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.132"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("E51").Value = "  +3.38%  "
